# Two independent fixes bundled in this commit:
#
# 1. On the "Repayment schedule" sheet, add the missing column O values
#    (mirroring the neighboring M/N "In Advance"/"Late" columns) so the
#    application's repayment-schedule validation has a complete row.
#
# 2. On the "Summary" sheet, the saved selection/view state pointed at a
#    stray range (A7:XFD15); reset it to the single cell C3.

$wb = $excel.ActiveWorkbook

# Remember which sheet/tab was active so we can restore it at the end -
# these selection/formatting tweaks shouldn't change which tab is shown.
$originalActiveSheet = $wb.ActiveSheet.Name

# --- Fix 1: Repayment schedule column O -------------------------------
$ws = $wb.Worksheets.Item("Repayment schedule")

# Pick up the same cell formatting used by column N (wrap text, vertically
# centered, "General" number format -> style index 10) for the new column O
# cells, without minting a brand-new/unused style entry.
$ws.Range("N2:N8").Copy()
$ws.Range("O2:O8").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 2 (the disbursement row) stays blank in column O; rows 3-8 (the actual
# repayment installments) get a 0 value, same as columns M and N.
$ws.Range("O3").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("O8").Value = 0

# --- Fix 2: Summary sheet selection -----------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("C3").Select()

# Restore the tab that was active before these edits.
$wb.Worksheets.Item($originalActiveSheet).Activate()
